$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("notebook")
Write-Host $ws.Name
Write-Host $wb.Worksheets.Count
